$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they stay text (matches source inlineStr)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values
$ws.Range('D2').Value = '62.500.76'
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').Value = '3.139.15'
$ws.Range('E3').Value = '  -5.20%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '586.39'
$ws.Range('E5').Value = '  -2.75%  '
$ws.Range('D6').Value = '135.38'
$ws.Range('E6').Value = '  -4.98%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.133.28'
$ws.Range('E8').Value = '  -5.43%  '
$ws.Range('E9').Value = '  -2.56%  '
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  -5.99%  '
$ws.Range('D11').Value = '5.23'
$ws.Range('E11').Value = '  -4.64%  '
$ws.Range('D12').Value = '0.452'
$ws.Range('E12').Value = '  -3.95%  '
$ws.Range('D13').Value = '0.0000233'
$ws.Range('E13').Value = '  -6.01%  '
$ws.Range('D14').Value = '33.78'
$ws.Range('E14').Value = '  -2.72%  '
$ws.Range('D15').Value = '3.639.75'
$ws.Range('E15').Value = '  -5.65%  '
$ws.Range('D16').Value = '0.118'
$ws.Range('E16').Value = '  -1.86%  '
$ws.Range('D17').Value = '3.152.21'
$ws.Range('E17').Value = '  -4.78%  '
$ws.Range('D18').Value = '62.465.02'
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('D19').Value = '6.53'
$ws.Range('E19').Value = '  -4.95%  '
$ws.Range('D20').Value = '451.86'
$ws.Range('E20').Value = '  -5.97%  '
$ws.Range('D21').Value = '13.88'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').Value = '0.696'
$ws.Range('E22').Value = '  -4.93%  '
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  -6.33%  '
$ws.Range('D24').Value = '13.29'
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('D25').Value = '83.24'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').Value = '2.67'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').Value = '7.66'
$ws.Range('E29').Value = '  -5.81%  '
$ws.Range('D30').Value = '6.69'
$ws.Range('E30').Value = '  -8.29%  '
$ws.Range('D31').Value = '2.00'
$ws.Range('E31').Value = '  -7.68%  '
$ws.Range('D32').Value = '27.06'
$ws.Range('E32').Value = '  -5.85%  '
$ws.Range('E33').Value = '  -3.59%  '
$ws.Range('D34').Value = '2.36'
$ws.Range('E34').Value = '  -7.41%  '
$ws.Range('D35').Value = '1.02'
$ws.Range('E35').Value = '  -7.54%  '
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  -4.08%  '
$ws.Range('D37').Value = '51.04'
$ws.Range('E37').Value = '  -4.34%  '
$ws.Range('D38').Value = '0.0₃0699'
$ws.Range('E38').Value = '  -6.23%  '
$ws.Range('D39').Value = '0.0384'
$ws.Range('E39').Value = '  -4.20%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '2.67'
$ws.Range('E40').Value = '  -3.16%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '395.57'
$ws.Range('E41').Value = '  -8.90%  '
$ws.Range('D42').Value = '8.00'
$ws.Range('E42').Value = '  -4.24%  '
$ws.Range('E43').Value = '  -3.23%  '
$ws.Range('D44').Value = '2.742.46'
$ws.Range('E44').Value = '  -10.87%  '
$ws.Range('E45').Value = '  -6.10%  '
$ws.Range('E47').Value = '  -4.52%  '
$ws.Range('D48').Value = '124.87'
$ws.Range('E48').Value = '  -3.70%  '
$ws.Range('D49').Value = '25.09'
$ws.Range('E49').Value = '  -4.90%  '
$ws.Range('D50').Value = '34.32'
$ws.Range('E50').Value = '  -5.03%  '
$ws.Range('D51').Value = '0.110'
$ws.Range('E51').Value = '  -3.85%  '
